# "Timesheet Changes by Ruchika" - fills in the remaining week's hours for
# rows 28-31 of the "February 2013" sheet, which belong to Ruchika Sharma
# (MT2012119): columns P, R, S, T, U get numbers, column Q (Saturday,
# 11-Feb-2012) gets the "OFF" label like the other OFF columns in the sheet,
# and column V is reformatted to the plain bordered style used elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("February 2013")

# --- Row 28 (Activity P02) --------------------------------------------
$ws.Range("P28").Value = 0
$ws.Range("R28").Value = 0
$ws.Range("S28").Value = 0
$ws.Range("T28").Value = 0
$ws.Range("U28").Value = 0

# --- Row 29 (Activity P03) --------------------------------------------
$ws.Range("P29").Value = 6
$ws.Range("R29").Value = 0
$ws.Range("S29").Value = 2
$ws.Range("T29").Value = 0
$ws.Range("U29").Value = 0

# --- Row 30 (Activity P04) --------------------------------------------
$ws.Range("P30").Value = 0
$ws.Range("R30").Value = 2
$ws.Range("S30").Value = 0
$ws.Range("T30").Value = 3
$ws.Range("U30").Value = 3

# --- Row 31 (Meetings) --------------------------------------------------
$ws.Range("P31").Value = 1
$ws.Range("R31").Value = 1.3
$ws.Range("S31").Value = 0
$ws.Range("T31").Value = 1.4
$ws.Range("U31").Value = 2

# Column Q (the Saturday "OFF" day) - copy the shaded "OFF" format already
# used further down the sheet (e.g. Q36) so the new cells pick up the same
# style, then write the label.
$ws.Range("Q36").Copy() | Out-Null
$ws.Range("Q28:Q31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("Q28").Value = "OFF"
$ws.Range("Q29").Value = "OFF"
$ws.Range("Q30").Value = "OFF"
$ws.Range("Q31").Value = "OFF"

# Column V for these rows switches to the plain bordered style used by the
# rest of the "OFF" block further down (e.g. V43) - format only, no value.
$ws.Range("V43").Copy() | Out-Null
$ws.Range("V28:V31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the view: scrolled down/right a bit further and the active
# selection moved from U38 to V30.
$ws.Range("V30").Select()
